$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.133.04'
$ws.Range('E2').Value = '  -10.91%  '
$ws.Range('D3').Value = '2.270.70'
$ws.Range('E3').Value = '  -21.92%  '
$ws.Range('D4').Value = '''0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '''449.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -15.09%  '
$ws.Range('D6').Value = '''128.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -11.32%  '
$ws.Range('D7').Value = '''0.996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = '''0.472'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -14.96%  '
$ws.Range('D9').Value = '2.252.15'
$ws.Range('E9').Value = '  -22.81%  '
$ws.Range('D10').Value = '''5.35'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -10.99%  '
$ws.Range('D11').Value = '''0.0919'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -15.46%  '
$ws.Range('E12').Value = '  -15.03%  '
$ws.Range('E13').Value = '  -3.12%  '
$ws.Range('D14').Value = '2.643.43'
$ws.Range('E14').Value = '  -22.62%  '
$ws.Range('D15').Value = '54.193.75'
$ws.Range('E15').Value = '  -10.71%  '
$ws.Range('D16').Value = '''18.83'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -17.45%  '
$ws.Range('E17').Value = '  -15.52%  '
$ws.Range('D18').Value = '2.267.56'
$ws.Range('E18').Value = '  -22.21%  '
$ws.Range('E19').Value = '  -19.48%  '
$ws.Range('D20').Value = '''302.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -16.34%  '
$ws.Range('D21').Value = '''9.46'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -19.15%  '
$ws.Range('D22').Value = '''0.999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').Value = '''5.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('D24').Value = '''5.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -20.05%  '
$ws.Range('D25').Value = '''55.69'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -14.03%  '
$ws.Range('D26').Value = '''0.974'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.50%  '
$ws.Range('D27').Value = '''0.158'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -13.36%  '
$ws.Range('D28').Value = '''0.373'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -18.10%  '
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').Value = '''0.996'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '''6.81'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -13.16%  '
$ws.Range('D31').Value = '0.0₃0707'
$ws.Range('E31').Value = '  -18.40%  '
$ws.Range('D32').Value = '''144.49'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.15%  '
$ws.Range('D33').Value = '''16.87'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -14.65%  '
$ws.Range('D34').Value = '''1.36'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -19.46%  '
$ws.Range('D35').Value = '''4.73'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -15.36%  '
$ws.Range('D36').Value = '''3.63'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -17.74%  '
$ws.Range('D37').Value = '''0.841'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -16.94%  '
$ws.Range('E38').Value = '  -17.14%  '
$ws.Range('D39').Value = '''0.990'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.70%  '
$ws.Range('D40').Value = '''32.95'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -12.67%  '
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('E42').Value = '  -16.36%  '
$ws.Range('D43').Value = '''3.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -15.83%  '
$ws.Range('D44').Value = '1.927.85'
$ws.Range('E44').Value = '  -15.79%  '
$ws.Range('D45').Value = '''0.0497'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -14.98%  '
$ws.Range('E46').Value = '  -13.34%  '
$ws.Range('D47').Value = '''0.505'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -22.38%  '
$ws.Range('D48').Value = '''0.0809'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -12.73%  '
$ws.Range('D49').Value = '''16.16'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -22.15%  '
$ws.Range('D50').Value = '''4.04'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -19.85%  '
$ws.Range('E51').Value = '  -3.21%  '
